# Delete the record for Kayıt No 11362186 from both the master "Kayitlar"
# sheet and the filtered "Merkez İlçe" sheet, shifting the following rows up.

$wb = $excel.ActiveWorkbook

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$wsMerkez = $wb.Worksheets.Item("Merkez İlçe")

$wsKayitlar.Rows.Item(1318).Delete()
$wsMerkez.Rows.Item(779).Delete()
